$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matching source formatting)
$dCells = @("D2", "D3", "D5", "D6", "D10", "D11", "D12", "D13", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D27", "D29", "D31", "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D42", "D43", "D44", "D46", "D48", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.118.71"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "3.499.01"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "597.71"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "175.20"
$ws.Range("E6").Value = "  +3.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -0.78%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "7.15"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("D11").Value = "0.430"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "4.108.19"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "30.85"
$ws.Range("E13").Value = "  +8.86%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "67.117.84"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  -1.85%  "
$ws.Range("D17").Value = "3.491.67"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "6.28"
$ws.Range("E18").Value = "  -0.73%  "
$ws.Range("D19").Value = "14.45"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").Value = "393.07"
$ws.Range("E20").Value = "  -0.73%  "
$ws.Range("D21").Value = "7.98"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "73.26"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").Value = "10.17"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -2.68%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").Value = "23.62"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").Value = "7.36"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("D36").Value = "163.14"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "0.879"
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("D40").Value = "27.47"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("D41").Value = "4.65"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").Value = "0.0730"
$ws.Range("E42").Value = "  -1.92%  "
$ws.Range("D43").Value = "26.06"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").Value = "2.795.36"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  -1.12%  "
$ws.Range("D46").Value = "42.42"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("D48").Value = "339.81"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "33.45"
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").Value = "0.843"
$ws.Range("E51").Value = "  -1.87%  "

Write-Output "Applied cryptos update"
